$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Started" value from D6 (row 6 now only has B6/C6)
$ws.Range("D6").ClearContents() | Out-Null

# New row 12: "Hacer un log" / "De lo que se le envía al servidor" / "ok"
$ws.Range("B12").Value = "Hacer un log"
$ws.Range("C12").Value = "De lo que se le envía al servidor"
$ws.Range("D12").Value = "ok"

# New row 13: D13 ("Start") set before B13 ("Despacahar las urls") so the
# shared-string table ends up in the same order as the target workbook.
$ws.Range("D13").Value = "Start"
$ws.Range("B13").Value = "Despacahar las urls"

# Column layout: add a narrow column A, and widen column C (drop its bestFit autosize).
# The host snaps ColumnWidth to an internal 1/6-character pixel grid, so these inputs
# are chosen to land on the grid points closest to the target widths (~5.71 and ~67.29).
$ws.Columns.Item(1).ColumnWidth = 4.833333333333334
$ws.Columns.Item(3).ColumnWidth = 66.5

# Update the active selection to D14, matching the post-edit cursor position
$ws.Range("D14").Select() | Out-Null

$wb.Save()
